$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.425.15'
$ws.Range("E2").Value = '  +2.19%  '

# Row 3
$ws.Range("D3").Value = '1.797.33'
$ws.Range("E3").Value = '  +2.85%  '

# Row 4
$ws.Range("E4").Value = '  +0.50%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '337.61'
$ws.Range("E5").Value = '  +0.52%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.19%  '

# Row 7
$ws.Range("E7").Value = '  +1.33%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3456'

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.60'
$ws.Range("E9").Value = '  +0.64%  '

# Row 10
$ws.Range("E10").Value = '  +0.73%  '

# Row 11
$ws.Range("E11").Value = '  -0.07%  '

# Row 12
$ws.Range("E12").Value = '  +0.44%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.11'
$ws.Range("E13").Value = '  +7.30%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.473'
$ws.Range("E14").Value = '  +0.65%  '

# Row 15
$ws.Range("D15").Value = '1.795.96'
$ws.Range("E15").Value = '  +3.01%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.075'
$ws.Range("E16").Value = '  -0.19%  '

# Row 17
$ws.Range("E17").Value = '  +1.63%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06656'
$ws.Range("E18").Value = '  -1.28%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '84.71'
$ws.Range("E19").Value = '  +2.14%  '

# Row 20
$ws.Range("E20").Value = '  +0.19%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.511'
$ws.Range("E21").Value = '  +4.38%  '

# Row 22
$ws.Range("E22").Value = '  +3.02%  '

# Row 23
$ws.Range("D23").Value = '27.400.41'
$ws.Range("E23").Value = '  +2.28%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.53'
$ws.Range("E24").Value = '  -2.37%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.430'
$ws.Range("E25").Value = '  -1.51%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.506'
$ws.Range("E26").Value = '  +2.15%  '

# Row 27
$ws.Range("E27").Value = '  +4.59%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.39'
$ws.Range("E28").Value = '  +8.62%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '152.20'
$ws.Range("E29").Value = '  +0.10%  '

# Row 30
$ws.Range("D30").Value = '2.000.43'
$ws.Range("E30").Value = '  +3.24%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '133.79'
$ws.Range("E31").Value = '  +0.83%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.056'
$ws.Range("E32").Value = '  -1.48%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.124'
$ws.Range("E33").Value = '  +0.68%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08668'
$ws.Range("E34").Value = '  +0.56%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.25'
$ws.Range("E35").Value = '  +2.21%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.657'
$ws.Range("E36").Value = '  -2.13%  '

# Row 37
$ws.Range("B37").Value = 'TheSandbox'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6906'
$ws.Range("E37").Value = '  +9.80%  '

# Row 38
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.447'
$ws.Range("E38").Value = '  -0.35%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.866'
$ws.Range("E39").Value = '  +4.15%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06367'
$ws.Range("E40").Value = '  +1.28%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2205'
$ws.Range("E41").Value = '  +0.83%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.02333'
$ws.Range("E42").Value = '  -1.15%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.279'
$ws.Range("E43").Value = '  +4.22%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.41'
$ws.Range("E44").Value = '  +0.24%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6439'
$ws.Range("E45").Value = '  +5.40%  '

# Row 46
$ws.Range("E46").Value = '  +0.28%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.867'
$ws.Range("E47").Value = '  -1.89%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.129'
$ws.Range("E48").Value = '  +2.41%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '130.34'
$ws.Range("E49").Value = '  +0.62%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07206'
$ws.Range("E50").Value = '  -0.44%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.82'
$ws.Range("E51").Value = '  +2.30%  '
